$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 31   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/2/2024  Through  12/8/2024"

# --- Numeric / percentage value updates across the crime stats table (rows 14-33) ---
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = -100
$ws.Range("J14").Value = 8
$ws.Range("K14").Value = -25
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -30.769230769230
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -62.5
$ws.Range("I16").Value = 148
$ws.Range("J16").Value = 164
$ws.Range("K16").Value = -9.756097560975
$ws.Range("L16").Value = -21.693121693121
$ws.Range("M16").Value = -35.930735930735
$ws.Range("N16").Value = 35.779816513761
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 23.076923076923
$ws.Range("I17").Value = 220
$ws.Range("J17").Value = 233
$ws.Range("K17").Value = -5.579399141630
$ws.Range("L17").Value = -13.385826771653
$ws.Range("M17").Value = 29.411764705882
$ws.Range("N17").Value = 147.191011235955
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 107
$ws.Range("J18").Value = 109
$ws.Range("K18").Value = -1.834862385321
$ws.Range("L18").Value = 8.080808080808
$ws.Range("M18").Value = 5.940594059405
$ws.Range("N18").Value = -11.570247933884
$ws.Range("C19").Value = 7
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = -3.703703703703
$ws.Range("I19").Value = 329
$ws.Range("J19").Value = 342
$ws.Range("K19").Value = -3.801169590643
$ws.Range("L19").Value = 17.081850533807
$ws.Range("M19").Value = 30.039525691699
$ws.Range("N19").Value = 327.272727272727
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 102
$ws.Range("J20").Value = 135
$ws.Range("K20").Value = -24.444444444444
$ws.Range("L20").Value = -21.538461538461
$ws.Range("M20").Value = 59.375
$ws.Range("N20").Value = 3.030303030303
$ws.Range("C21").Value = 12
$ws.Range("E21").Value = -36.842105263157
$ws.Range("F21").Value = 59
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = -19.178082191780
$ws.Range("I21").Value = 921
$ws.Range("J21").Value = 998
$ws.Range("K21").Value = -7.715430861723
$ws.Range("L21").Value = -4.855371900826
$ws.Range("M21").Value = 9.512485136741
$ws.Range("N21").Value = 81.299212598425
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -88.888888888888
$ws.Range("I22").Value = 16
$ws.Range("K22").Value = -44.827586206896
$ws.Range("L22").Value = -40.740740740740
$ws.Range("M22").Value = 14.285714285714
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 2
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = 80
$ws.Range("L23").Value = 35
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = -40
$ws.Range("F24").Value = 79
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = 51.923076923076
$ws.Range("I24").Value = 867
$ws.Range("J24").Value = 867
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = -29.569455727051
$ws.Range("M24").Value = 56.781193490054
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 184.615384615385
$ws.Range("I25").Value = 294
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 47
$ws.Range("L25").Value = -54.629629629629
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 16.666666666666
$ws.Range("I26").Value = 328
$ws.Range("J26").Value = 348
$ws.Range("K26").Value = -5.747126436781
$ws.Range("L26").Value = -12.765957446808
$ws.Range("M26").Value = -27.593818984547
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 0
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 35
$ws.Range("J28").Value = 41
$ws.Range("K28").Value = -14.634146341463
$ws.Range("L28").Value = -25.531914893617
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 13
$ws.Range("K29").Value = -53.846153846153
$ws.Range("L29").Value = -45.454545454545
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 9
$ws.Range("K30").Value = -55.555555555555
$ws.Range("L30").Value = -60
$ws.Range("L33").Value = 100

# --- Cells switching from numeric back to the "no data" text placeholders ---
# Donor cells C14 (text "0", style 13) and C15 style... use stable untouched style-13 donor cells
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4122)

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)

$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
